$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Priority" (row 59) and "Registered" (row 60) service rows.
# This shifts "International Tracked Packet" up to row 59 and the
# footnote rows up to rows 60-63.
$ws.Rows("59:60").Delete()

# Column B ("Customer number type") for the "Cargo Norway international"
# row should read "Main customer number" instead of repeating the
# service-family name.
$ws.Range("B56").Value = "Main customer number"

# Narrow column B now that the longest entries were removed.
$ws.Columns("B").ColumnWidth = 21.365885416666668

# Keep the autofilter and the hidden _FilterDatabase defined name in sync
# with the new, smaller data range (2 fewer rows). The sheet already has
# an active AutoFilter, so turn it off first - otherwise re-applying it on
# the same range just toggles filtering off instead of resizing it.
$ws.AutoFilterMode = $false
$ws.Range("A1:P64").AutoFilter()
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Booking & SG API'!`$A`$1:`$P`$64"
